# Updated cryptos list: refresh Price (column D) and Volume(1h) (column E)
# values for the coinranking.com snapshot rows (rows 2-51).
# For D-column updates whose new value parses as a plain number (e.g. "581.92"),
# force the cell to Text format first so Excel keeps the original decimal
# formatting (e.g. trailing zeros) instead of silently converting it to a
# numeric value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.882.60'
$ws.Range("E2").Value = '  -0.47%  '
$ws.Range("D3").Value = '3.258.79'
$ws.Range("E3").Value = '  -0.60%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '581.92'
$ws.Range("E5").Value = '  -0.88%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '185.53'
$ws.Range("E6").Value = '  +0.14%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.603'
$ws.Range("E8").Value = '  +0.44%  '
$ws.Range("D9").Value = '3.256.52'
$ws.Range("E9").Value = '  -0.61%  '
$ws.Range("E10").Value = '  -3.09%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.59'
$ws.Range("E11").Value = '  -2.04%  '
$ws.Range("E12").Value = '  -1.55%  '
$ws.Range("D13").Value = '3.822.08'
$ws.Range("E13").Value = '  -0.68%  '
$ws.Range("E14").Value = '  -0.12%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.67'
$ws.Range("E15").Value = '  -3.56%  '
$ws.Range("D16").Value = '67.817.09'
$ws.Range("E16").Value = '  -0.55%  '
$ws.Range("E17").Value = '  -1.74%  '
$ws.Range("D18").Value = '3.250.23'
$ws.Range("E18").Value = '  -0.79%  '
$ws.Range("E19").Value = '  -2.24%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.59'
$ws.Range("E20").Value = '  -0.40%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '394.10'
$ws.Range("E21").Value = '  +3.06%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.60'
$ws.Range("E22").Value = '  -2.28%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '71.75'
$ws.Range("E23").Value = '  +0.51%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.999'
$ws.Range("E24").Value = '  -0.11%  '
$ws.Range("E25").Value = '  +0.10%  '
$ws.Range("E26").Value = '  -2.24%  '
$ws.Range("E27").Value = '  -3.12%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.56'
$ws.Range("E28").Value = '  -2.27%  '
$ws.Range("E29").Value = '  +1.43%  '
$ws.Range("E30").Value = '  -1.77%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.53'
$ws.Range("E31").Value = '  -5.01%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '22.68'
$ws.Range("E32").Value = '  -1.37%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.00'
$ws.Range("E33").Value = '  -3.09%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.25'
$ws.Range("E34").Value = '  -2.68%  '
$ws.Range("E35").Value = '  +0.03%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '162.91'
$ws.Range("E36").Value = '  -0.35%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.47'
$ws.Range("E37").Value = '  -4.26%  '
$ws.Range("E38").Value = '  +2.08%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '26.66'
$ws.Range("E39").Value = '  -0.10%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.810'
$ws.Range("E40").Value = '  -3.56%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.54'
$ws.Range("E41").Value = '  -1.91%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.44'
$ws.Range("E42").Value = '  -4.31%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0688'
$ws.Range("E43").Value = '  -0.65%  '
$ws.Range("E44").Value = '  -7.14%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '40.60'
$ws.Range("E45").Value = '  -1.58%  '
$ws.Range("D46").Value = '2.610.33'
$ws.Range("E46").Value = '  -0.72%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '24.83'
$ws.Range("E47").Value = '  -3.26%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '333.69'
$ws.Range("E48").Value = '  -2.78%  '
$ws.Range("E49").Value = '  -2.91%  '
$ws.Range("E51").Value = '  -0.95%  '
